# Update "想去人数" (interest count) values across the sheets, reflecting
# the data refresh recorded in the commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 15687
$ws1.Range("F9").Value  = 15433
$ws1.Range("F10").Value = 54
$ws1.Range("F11").Value = 9033
$ws1.Range("F12").Value = 385
$ws1.Range("F15").Value = 94
$ws1.Range("F20").Value = 55
$ws1.Range("F29").Value = 87
$ws1.Range("F32").Value = 412
$ws1.Range("F38").Value = 118
$ws1.Range("F39").Value = 5559

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 68

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 15687
$ws4.Range("F9").Value  = 15433
$ws4.Range("F10").Value = 54
$ws4.Range("F11").Value = 9033
$ws4.Range("F12").Value = 385
$ws4.Range("F15").Value = 94
$ws4.Range("F20").Value = 55
$ws4.Range("F29").Value = 87
$ws4.Range("F32").Value = 68
$ws4.Range("F34").Value = 412
$ws4.Range("F40").Value = 118
$ws4.Range("F41").Value = 5559
